$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 previously had an (empty) Proxy value in column F; clear it so the
# cell no longer holds any content.
$ws.Cells.Item(42, 6).Value = ""

# New account rows appended to the log (A:Name, B:Email, C:Password,
# D:DOB, E:Gender, F:Proxy, G:Created At)
$rows = @(
    @("Michael Smith", "bapije3413@agiuse.com", "Pass8934", "15 September 1983", "female", "", "2025-04-19 23:46:41"),
    @("Michael Smith", "bapije3413@agiuse.com", "Pass8934", "15 September 1983", "female", "", "2025-04-19 23:46:45"),
    @("John Garcia", "johngarcia52@maildrop.cc", "Pass7741", "20 December 1989", "male", "", "2025-04-20 00:49:42"),
    @("John Garcia", "johngarcia52@maildrop.cc", "Pass7741", "20 December 1989", "male", "", "2025-04-20 00:49:42"),
    @("Emma Smith", "emmasmith62@maildrop.cc", "Pass6692", "20 June 1991", "male", "", "2025-04-20 00:50:56"),
    @("Sarah Brown", "oceubme753@1secmail.website", "Pass9123", "23 April 1998", "female", "", "2025-04-20 01:18:13"),
    @("Sarah Brown", "oceubme753@1secmail.website", "Pass9123", "23 April 1998", "female", "", "2025-04-20 01:18:13")
)

$startRow = 43
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[6]
}

# Last appended row (49) also carries an explicit, empty Proxy cell.
$ws.Cells.Item(49, 6).Value = ""
